$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Gaussian" section: rows 37-41 -----------------------------------
# Row 37: section headers (merged A37:D37 and F37:I37), centered + wrap text
$ws.Range("F37").Value2 = "Avg MFCC (13 coeff) + Delta + GMM [1] + 65 knn"
$ws.Range("A37").Value2 = "Avg MFCC (13 coeff) + Delta + Delta Delta + GMM (1) + 100 knn"

$ws.Range("A37:D37").HorizontalAlignment = -4108
$ws.Range("A37:D37").WrapText = $true
$ws.Range("F37:I37").HorizontalAlignment = -4108
$ws.Range("F37:I37").WrapText = $true
$ws.Range("A37:D37").Merge()
$ws.Range("F37:I37").Merge()
$ws.Rows.Item(37).RowHeight = 36.75

# Row 38: "EER" sub-headers (merged B38:D38 and G38:I38), centered
$ws.Range("B38").Value2 = "EER"
$ws.Range("G38").Value2 = "EER"
$ws.Range("B38:D38").HorizontalAlignment = -4108
$ws.Range("G38:I38").HorizontalAlignment = -4108
$ws.Range("B38:D38").Merge()
$ws.Range("G38:I38").Merge()

# Row 39: column headers
$ws.Range("A39").Value2 = "Train"
$ws.Range("B39").Value2 = "Test: Read"
$ws.Range("C39").Value2 = "Test: Phone"
$ws.Range("D39").Value2 = "Test: Mismatch"
$ws.Range("F39").Value2 = "Train"
$ws.Range("G39").Value2 = "Test: Read"
$ws.Range("H39").Value2 = "Test: Phone"
$ws.Range("I39").Value2 = "Test: Mismatch"

# Row 40: "Read" data
$ws.Range("A40").Value2 = "Read"
$ws.Range("B40").Value2 = 22.222
$ws.Range("C40").Value2 = 29.649
$ws.Range("D40").Value2 = 37.777
$ws.Range("F40").Value2 = "Read"
$ws.Range("G40").Value2 = 15.55
$ws.Range("H40").Value2 = 29.41
$ws.Range("I40").Value2 = 35.55

# Row 41: "Phone" data
$ws.Range("A41").Value2 = "Phone"
$ws.Range("B41").Value2 = 33.333
$ws.Range("C41").Value2 = 25
$ws.Range("D41").Value2 = 46.654
$ws.Range("F41").Value2 = "Phone"
$ws.Range("G41").Value2 = 34.3
$ws.Range("H41").Value2 = 23.33
$ws.Range("I41").Value2 = 47.65

# --- View state: scroll position + selection -------------------------------
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("D41").Select()
